# Add a new worksheet "Autonomous_temporary" right after "Autonomous",
# containing an auto-generated export of DBC signal tables.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Autonomous")

$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "Autonomous_temporary"

# Match outline defaults used on the sheet (summary rows below / summary
# columns to the right of the detail, Excel's usual default).
$new.Outline.SummaryRow = 1
$new.Outline.SummaryColumn = 1

# Match column widths used on the new sheet (xml width = ColumnWidth + 0.83)
$widths = @(21,11,15,12,8,8,8,5,5,6,9)
for ($i = 1; $i -le 11; $i++) {
    $new.Columns.Item($i).ColumnWidth = $widths[$i - 1] - 0.83
}

# Messages, in order, as they appear in the generated sheet.
# Each entry: MessageName, MessageId, SignalName, LengthInBits
$messages = @(
    @("ACU_MS",     "0x51",  "Mission_Select", 8),
    @("JETSON_MS",  "0x61",  "Mission_Select", 8),
    @("VCU_RPM",    "0x510", "RPM",            16),
    @("Target_RPM", "0x500", "RPM",            16),
    @("ACU_IGN",    "0x71",  "IGN",            8),
    @("RD_JETSON",  "0x512", "RD",             8),
    @("AS_STATE",   "0x502", "STATE",          8),
    @("VCU_HV",     "0x81",  "HV",             8)
)

$row = 1
foreach ($m in $messages) {
    $msgName = $m[0]
    $msgId = $m[1]
    $sigName = $m[2]
    $length = $m[3]

    $headerRow = $row
    $colHeaderRow = $row + 1
    $dataRow = $row + 2

    # Paste formatting for this 3-row block (header / column-header / data).
    # Reuse the existing header/sub-header/data-row formatting (styles
    # 1,2,3) from the "Autonomous" sheet. The message/id row only spans
    # columns A:B, while the column-header and data rows span the full A:K.
    $ws1.Range("A1:B1").Copy()
    $new.Range("A$headerRow`:B$headerRow").PasteSpecial(-4122)
    $ws1.Range("A2:K3").Copy()
    $new.Range("A$colHeaderRow`:K$dataRow").PasteSpecial(-4122)

    # Row 1 of block: message name / id
    $new.Range("A$headerRow").Value = "Message: $msgName"
    $new.Range("B$headerRow").Value = "ID: $msgId"

    # Row 2 of block: column headers
    $new.Range("A$colHeaderRow").Value = "Signal Name"
    $new.Range("B$colHeaderRow").Value = "Start Bit"
    $new.Range("C$colHeaderRow").Value = "Length (bits)"
    $new.Range("D$colHeaderRow").Value = "Byte Order"
    $new.Range("E$colHeaderRow").Value = "Signed"
    $new.Range("F$colHeaderRow").Value = "Factor"
    $new.Range("G$colHeaderRow").Value = "Offset"
    $new.Range("H$colHeaderRow").Value = "Min"
    $new.Range("I$colHeaderRow").Value = "Max"
    $new.Range("J$colHeaderRow").Value = "Unit"
    $new.Range("K$colHeaderRow").Value = "Choices"

    # Row 3 of block: the single signal's data
    $new.Range("A$dataRow").Value = $sigName
    $new.Range("B$dataRow").Value = 0
    $new.Range("C$dataRow").Value = $length
    $new.Range("D$dataRow").Value = "Intel"
    $new.Range("E$dataRow").Value = $false
    $new.Range("F$dataRow").Value = 1
    $new.Range("G$dataRow").Value = 0

    $row = $row + 4
}

$excel.CutCopyMode = 0

# Keep "Autonomous" as the active sheet/tab, matching the original workbook.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
